$wb = $excel.ActiveWorkbook

# This script updates market-price derived columns (H:N) on each class sheet
# (currentAveragePrice, currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ)
# to reflect a refreshed data pull, per the scheduled-runner commit.

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 94.083336
$ws.Range("I2").Value = 84.454544
$ws.Range("J2").Value = 200
$ws.Range("K2").Value = 84.454544
$ws.Range("L2").Value = 200
$ws.Range("M2").Value = 28.545456
$ws.Range("N2").Value = -426

$ws.Range("H4").Value = 209.5
$ws.Range("I4").Value = 185
$ws.Range("J4").Value = 266.66666
$ws.Range("K4").Value = 185
$ws.Range("L4").Value = 266.66666
$ws.Range("M4").Value = -71
$ws.Range("N4").Value = -494.66666

$ws.Range("H9").Value = 68.8
$ws.Range("I9").Value = 69.5
$ws.Range("K9").Value = 69.5
$ws.Range("M9").Value = 99.5

$ws.Range("H12").Value = 290.7647
$ws.Range("I12").Value = 148.09091
$ws.Range("J12").Value = 4999
$ws.Range("K12").Value = 148.09091
$ws.Range("L12").Value = 4999
$ws.Range("M12").Value = 21.90908999999999
$ws.Range("N12").Value = -5339

$ws.Range("H17").Value = 5884978
$ws.Range("J17").Value = 6252708.5
$ws.Range("L17").Value = 18758125.5
$ws.Range("N17").Value = -18758461.5

$ws.Range("H18").Value = 5322
$ws.Range("J18").Value = 450
$ws.Range("L18").Value = 450
$ws.Range("N18").Value = -1018

$ws.Range("H70").Value = 3528.5715
$ws.Range("J70").Value = 4675
$ws.Range("L70").Value = 14025
$ws.Range("N70").Value = -14565

$ws.Range("H73").Value = 3528.5715
$ws.Range("J73").Value = 4675
$ws.Range("L73").Value = 14025
$ws.Range("N73").Value = -15897

$ws.Range("H87").Value = 61759.332
$ws.Range("J87").Value = 61759.332
$ws.Range("L87").Value = 61759.332
$ws.Range("N87").Value = -64255.332

$ws.Range("H90").Value = 61759.332
$ws.Range("J90").Value = 61759.332
$ws.Range("L90").Value = 185277.996
$ws.Range("N90").Value = -197757.996

$ws.Range("H106").Value = 3842.7693
$ws.Range("I106").Value = 3190.25
$ws.Range("K106").Value = 3190.25
$ws.Range("M106").Value = -2559.25

$ws.Range("H116").Value = 3381.2727
$ws.Range("I116").Value = 3142.1428
$ws.Range("J116").Value = 3799.75
$ws.Range("K116").Value = 3142.1428
$ws.Range("L116").Value = 3799.75
$ws.Range("M116").Value = 299.8571999999999
$ws.Range("N116").Value = -10683.75

$ws.Range("H125").Value = 763.5
$ws.Range("I125").Value = 870.25
$ws.Range("J125").Value = 692.3333
$ws.Range("K125").Value = 7832.25
$ws.Range("L125").Value = 6230.9997
$ws.Range("M125").Value = -5372.25
$ws.Range("N125").Value = -11150.9997

$ws.Range("H132").Value = 1486.2916
$ws.Range("I132").Value = 1085.159
$ws.Range("K132").Value = 3255.477
$ws.Range("M132").Value = -725.4770000000003

$ws.Range("H137").Value = 2902.25
$ws.Range("I137").Value = 2602.6428
$ws.Range("K137").Value = 7807.928400000001
$ws.Range("M137").Value = -5257.928400000001


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 6037.615
$ws.Range("J45").Value = 7498.1665
$ws.Range("L45").Value = 7498.1665
$ws.Range("N45").Value = -8252.1665

$ws.Range("H61").Value = 10862.929
$ws.Range("I61").Value = 10862.929
$ws.Range("K61").Value = 10862.929
$ws.Range("M61").Value = -10650.929

$ws.Range("H122").Value = 1746.5
$ws.Range("I122").Value = 1912
$ws.Range("J122").Value = 1498.25
$ws.Range("K122").Value = 5736
$ws.Range("L122").Value = 4494.75
$ws.Range("M122").Value = -3286
$ws.Range("N122").Value = -9394.75

$ws.Range("H135").Value = 65590.5
$ws.Range("J135").Value = 65590.5
$ws.Range("L135").Value = 65590.5
$ws.Range("N135").Value = -75730.5

$ws.Range("H136").Value = 10862.929
$ws.Range("I136").Value = 10862.929
$ws.Range("K136").Value = 32588.787
$ws.Range("M136").Value = -30038.787


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3806.6155
$ws.Range("I20").Value = 3197.8667
$ws.Range("J20").Value = 4636.727
$ws.Range("K20").Value = 3197.8667
$ws.Range("L20").Value = 4636.727
$ws.Range("M20").Value = -2950.8667
$ws.Range("N20").Value = -5130.727

$ws.Range("H94").Value = 761.5
$ws.Range("I94").Value = 601.25
$ws.Range("K94").Value = 601.25
$ws.Range("M94").Value = -150.25

$ws.Range("H105").Value = 3346.7144
$ws.Range("I105").Value = 3537.4
$ws.Range("K105").Value = 3537.4
$ws.Range("M105").Value = -1790.4

$ws.Range("H107").Value = 1843.9025
$ws.Range("I107").Value = 1815
$ws.Range("J107").Value = 1933.5
$ws.Range("K107").Value = 1815
$ws.Range("L107").Value = 1933.5
$ws.Range("M107").Value = 105
$ws.Range("N107").Value = -5773.5

$ws.Range("H134").Value = 4761.3184
$ws.Range("I134").Value = 4262.45
$ws.Range("K134").Value = 12787.35
$ws.Range("M134").Value = -10252.35


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 62508580
$ws.Range("I31").Value = 111115980
$ws.Range("J31").Value = 13342.714
$ws.Range("K31").Value = 111115980
$ws.Range("L31").Value = 13342.714
$ws.Range("M31").Value = -111115685
$ws.Range("N31").Value = -13932.714

$ws.Range("H34").Value = 62508580
$ws.Range("I34").Value = 111115980
$ws.Range("J34").Value = 13342.714
$ws.Range("K34").Value = 111115980
$ws.Range("L34").Value = 13342.714
$ws.Range("M34").Value = -111115778
$ws.Range("N34").Value = -13746.714

$ws.Range("H94").Value = 1699.7858
$ws.Range("I94").Value = 1584.8
$ws.Range("J94").Value = 1763.6666
$ws.Range("K94").Value = 1584.8
$ws.Range("L94").Value = 1763.6666
$ws.Range("M94").Value = -1133.8
$ws.Range("N94").Value = -2665.6666

$ws.Range("H105").Value = 1405.4
$ws.Range("J105").Value = 1323.5
$ws.Range("L105").Value = 1323.5
$ws.Range("N105").Value = -4817.5

$ws.Range("H134").Value = 2915.3794
$ws.Range("I134").Value = 1779.64
$ws.Range("K134").Value = 5338.92
$ws.Range("M134").Value = -2803.92

$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("M135").ClearContents()

$ws.Range("H141").Value = 183713
$ws.Range("J141").Value = 212252.9
$ws.Range("L141").Value = 212252.9
$ws.Range("N141").Value = -222612.9


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 14286451
$ws.Range("I11").Value = 25001078
$ws.Range("J11").Value = 281.66666
$ws.Range("K11").Value = 75003234
$ws.Range("L11").Value = 844.9999799999999
$ws.Range("M11").Value = -75003094
$ws.Range("N11").Value = -1124.99998

$ws.Range("H34").Value = 424.83334
$ws.Range("J34").Value = 424.75
$ws.Range("L34").Value = 1274.25
$ws.Range("N34").Value = -1442.25

$ws.Range("H39").Value = 1343.75
$ws.Range("I39").Value = 1187.5
$ws.Range("K39").Value = 3562.5
$ws.Range("M39").Value = -3268.5

$ws.Range("H55").Value = 15502
$ws.Range("J55").Value = 30000
$ws.Range("L55").Value = 90000
$ws.Range("N55").Value = -90354

$ws.Range("H61").Value = 158.6
$ws.Range("I61").Value = 46.5
$ws.Range("K61").Value = 139.5
$ws.Range("M61").Value = 75.5

$ws.Range("H113").Value = 2430.889
$ws.Range("J113").Value = 2969.8572
$ws.Range("L113").Value = 8909.5716
$ws.Range("N113").Value = -13249.5716

$ws.Range("H121").Value = 8500698
$ws.Range("I121").Value = 125356
$ws.Range("K121").Value = 376068
$ws.Range("M121").Value = -374758

$ws.Range("H131").Value = 29414904
$ws.Range("I131").Value = 166667100
$ws.Range("J131").Value = 3718.0715
$ws.Range("K131").Value = 500001300
$ws.Range("L131").Value = 11154.2145
$ws.Range("M131").Value = -499996260
$ws.Range("N131").Value = -21234.2145


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 7247.857
$ws.Range("I102").Value = 5447
$ws.Range("K102").Value = 5447
$ws.Range("M102").Value = -3825

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2049.375
$ws.Range("I7").Value = 1519
$ws.Range("J7").Value = 2933.3333
$ws.Range("K7").Value = 1519
$ws.Range("L7").Value = 2933.3333
$ws.Range("M7").Value = -1407
$ws.Range("N7").Value = -3157.3333

$ws.Range("H16").Value = 658.8125
$ws.Range("I16").Value = 490.41666
$ws.Range("J16").Value = 1164
$ws.Range("K16").Value = 490.41666
$ws.Range("L16").Value = 1164
$ws.Range("M16").Value = -320.41666
$ws.Range("N16").Value = -1504

$ws.Range("H126").Value = 2049.375
$ws.Range("I126").Value = 1519
$ws.Range("J126").Value = 2933.3333
$ws.Range("K126").Value = 4557
$ws.Range("L126").Value = 8799.999899999999
$ws.Range("M126").Value = -2087
$ws.Range("N126").Value = -13739.9999

$ws.Range("H132").Value = 6209.069
$ws.Range("I132").Value = 6726.881
$ws.Range("J132").Value = 4849.8125
$ws.Range("K132").Value = 20180.643
$ws.Range("L132").Value = 14549.4375
$ws.Range("M132").Value = -17650.643
$ws.Range("N132").Value = -19609.4375

$ws.Range("H136").Value = 86959840
$ws.Range("I136").Value = 50002820
$ws.Range("K136").Value = 150008460
$ws.Range("M136").Value = -150005910


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1605.1482
$ws.Range("I107").Value = 1144.1578
$ws.Range("K107").Value = 3432.4734
$ws.Range("M107").Value = -1512.4734

$ws.Range("H126").Value = 2930.7568
$ws.Range("I126").Value = 2734.2666
$ws.Range("K126").Value = 8202.7998
$ws.Range("M126").Value = -5732.799800000001

$ws.Range("H132").Value = 2860.5386
$ws.Range("I132").Value = 2492.75
$ws.Range("J132").Value = 5149
$ws.Range("K132").Value = 7478.25
$ws.Range("L132").Value = 15447
$ws.Range("M132").Value = -4948.25
$ws.Range("N132").Value = -20507

